$d = $word.ActiveDocument
$sec = $d.Sections(1)
$f1 = $sec.Footers(1)
$result = $f1.Range.Find.Execute("2021", $true, $false, $false, $false, $false,
                   $true, 1, $false, "2022", 2)
Write-Host "Footer1 result=$result"
